$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 172, shifting existing rows 172:194 down to 173:195
$ws.Rows("172:172").Insert()

# Populate the newly inserted row 172 with a new Ciboulette price record
$ws.Range("A172").Value = 3
$ws.Range("B172").Value = "Femacal de La Calera"
$ws.Range("C172").Value = "Coquimbo"
$ws.Range("D172").Value = 44491
$ws.Range("E172").Value = 5
$ws.Range("F172").Value = 100112039
$ws.Range("G172").Value = "Ciboulette"
$ws.Range("H172").Value = "Sin especificar"
$ws.Range("I172").Value = "Primera"
$ws.Range("J172").Value = 160
$ws.Range("K172").Value = 1500
$ws.Range("L172").Value = 1500
$ws.Range("M172").Value = 1500
$ws.Range("N172").Value = "$/docena de atados"
$ws.Range("O172").Value = "Provincia de Quillota"
$ws.Range("P172").Value = 500
$ws.Range("Q172").Value = 3
$ws.Range("R172").Value = "Hortaliza"
